$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above the current row 121, pushing the existing
# weekly blocks (rows 121-180) down to rows 125-184.
$ws.Rows("121:124").Insert()

# New weekly block (date 44673) for "Agricola del Norte S.A. de Arica" /
# Pina / Caramelo, one row per quality grade.
$newRows = @(
    @{ Row = 121; Quality = "Especial"; Volumen = 200; PMin = 18000; PMax = 19000; PProm = 18500; Unidad = "$/caja 10 unidades"; PrecioKg = 1850; KgUnidad = 10 },
    @{ Row = 122; Quality = "Primera";  Volumen = 270; PMin = 18000; PMax = 19000; PProm = 18500; Unidad = "$/caja 12 unidades"; PrecioKg = 1542; KgUnidad = 12 },
    @{ Row = 123; Quality = "Segunda";  Volumen = 300; PMin = 18000; PMax = 19000; PProm = 18500; Unidad = "$/caja 14 unidades"; PrecioKg = 1321; KgUnidad = 14 },
    @{ Row = 124; Quality = "Tercera";  Volumen = 250; PMin = 18000; PMax = 19000; PProm = 18500; Unidad = "$/caja 16 unidades"; PrecioKg = 1156; KgUnidad = 16 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 1
    $ws.Cells.Item($row, 2).Value = "Agrícola del Norte S.A. de Arica"
    $ws.Cells.Item($row, 3).Value = "Arica y Parinacota"
    $ws.Cells.Item($row, 4).Value = 44673
    $ws.Cells.Item($row, 5).Value = 15
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100108
    $ws.Cells.Item($row, 8).Value = "Tropicales y subtropicales"
    $ws.Cells.Item($row, 9).Value = 100108005
    $ws.Cells.Item($row, 10).Value = "Piña"
    $ws.Cells.Item($row, 11).Value = "Caramelo"
    $ws.Cells.Item($row, 12).Value = $r.Quality
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = "Ecuador"
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = $r.KgUnidad
}
